# edit.ps1 - apply the Requirements.docx changes described by the diff.

function Find-ParagraphIndexByText($d, [string]$searchText) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text -like "*$searchText*") {
            return $i
        }
    }
    return -1
}

# Word-processingML namespace helper, used to build raw-XML fragments for
# Range.InsertXML so we get exact control over run/paragraph structure
# (tabs as <w:tab/>, lastRenderedPageBreak placement, etc.)
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. "5. Blocken stoppt Stamina Regeneration. ..." paragraph:
#    - add a left indent of 708 twips (35.4pt)
#    - drop the leading tab/space run
#    - replace the trailing "[-30% vom Angriff]" with the new sentence
# ---------------------------------------------------------------------------
$idx1 = Find-ParagraphIndexByText $d "Blocken stoppt Stamina Regeneration"
if ($idx1 -lt 0) { throw "Could not find the 'Blocken stoppt' paragraph" }
$p1 = $d.Paragraphs.Item($idx1)
$p1.Range.ParagraphFormat.LeftIndent = 35.4
$r1 = $p1.Range.Duplicate
$r1.MoveEnd(1, -1) # exclude the paragraph mark
$r1.Text = "5. Blocken stoppt Stamina Regeneration. Beim erfolgreichen Blocken wird von der Stamina der geblockte Damage abgezogen. Falls es keine Stamina gibt, werden Lebenspunkte abgezogen. "

# ---------------------------------------------------------------------------
# 2. "2.2.3 Gegner & Boss" heading: drop the stray <w:lastRenderedPageBreak/>.
#    A tiny round-trip text edit (append + revert marker, then restore)
#    forces the run(s) to be rewritten without the rendering artifact.
# ---------------------------------------------------------------------------
$idx2 = Find-ParagraphIndexByText $d "2.2.3 Gegner"
if ($idx2 -lt 0) { throw "Could not find the '2.2.3 Gegner' paragraph" }
$p2 = $d.Paragraphs.Item($idx2)
$r2 = $p2.Range.Duplicate
$r2.Find.Execute("2.2.3 Gegner", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r2.Text = "2.2.3 GegnerZZZ"
$r2b = $p2.Range.Duplicate
$r2b.Find.Execute("2.2.3 GegnerZZZ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r2b.Text = "2.2.3 Gegner"

# ---------------------------------------------------------------------------
# 3. Before "2.2.4.2 Kiste": insert a new "4. Stamina füllt sich zu 100%"
#    paragraph, and add a <w:lastRenderedPageBreak/> at the start of the
#    (now following) "2.2.4.2 Kiste" run.
# ---------------------------------------------------------------------------
$idx3 = Find-ParagraphIndexByText $d "2.2.4.2 Kiste"
if ($idx3 -lt 0) { throw "Could not find the '2.2.4.2 Kiste' paragraph" }
$p3 = $d.Paragraphs.Item($idx3)
$start3 = $p3.Range.Start
$rng3 = $d.Range($start3, $start3)
$uuml = [char]0x00FC
$xmlFrag3 = "<w:p $wns><w:r><w:tab/><w:t>4. Stamina f${uuml}llt sich zu 100%</w:t></w:r></w:p><w:p $wns/>"
$rng3.InsertXML($xmlFrag3)

# InsertXML leaves a trailing empty paragraph between the new paragraph and
# "2.2.4.2 Kiste" (an artifact of splitting at the very start of a range);
# remove it so Kiste directly follows the new paragraph, as in the diff.
$idxKiste = Find-ParagraphIndexByText $d "2.2.4.2 Kiste"
$idxBlank = $idxKiste - 1
$blankP = $d.Paragraphs.Item($idxBlank)
if ($blankP.Range.Text.Trim() -eq "") {
    $delRng = $d.Range($blankP.Range.Start, $blankP.Range.End)
    $delRng.Delete()
} else {
    throw "Expected a blank paragraph before 'Kiste', found: [$($blankP.Range.Text)]"
}

# Now add lastRenderedPageBreak to the start of the Kiste paragraph's bold run
$idx3b = Find-ParagraphIndexByText $d "2.2.4.2 Kiste"
$p3b = $d.Paragraphs.Item($idx3b)
$start3b = $p3b.Range.Start
$rng3b = $d.Range($start3b, $start3b)
$xmlFrag3b = "<w:p $wns><w:r><w:rPr><w:b/></w:rPr><w:lastRenderedPageBreak/></w:r></w:p>"
$rng3b.InsertXML($xmlFrag3b)

# ---------------------------------------------------------------------------
# 4. After "2.2.4.3 Level Design ... (Wird verfeinert)": insert the new
#    2.2.4.4 - 2.2.4.8 paragraphs (Platform / dyn. Platform / Spike /
#    dyn. Spike / Teleport) plus a trailing blank paragraph.
# ---------------------------------------------------------------------------
$idx4 = Find-ParagraphIndexByText $d "Wird verfeinert"
if ($idx4 -lt 0) { throw "Could not find the 'Wird verfeinert' paragraph" }
$p4 = $d.Paragraphs.Item($idx4 + 1)
$start4 = $p4.Range.Start
$rng4 = $d.Range($start4, $start4)

$ouml = [char]0x00F6
$szlig = [char]0x00DF

$xmlFrag4 = @"
<w:p $wns><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">2.2.4.4 </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Platform</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> -</w:t></w:r><w:r><w:t xml:space="preserve"> Der Spieler kann sich darauf bewegen und st${ouml}${szlig}t von jeder Seite dagegen.</w:t></w:r></w:p><w:p $wns><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">2.2.4.5 Dynamische </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Platform</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> -</w:t></w:r><w:r><w:t xml:space="preserve"> Wie eine </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Platform</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> nur das diese sich in einem festen Bereich bewegt</w:t></w:r></w:p><w:p $wns><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>2.2.4.6 Spike -</w:t></w:r><w:r><w:t xml:space="preserve"> Spieler kann durchlaufen und Schaden bekommen</w:t></w:r></w:p><w:p $wns><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>2.2.4.7 Dynamischer Spike -</w:t></w:r><w:r><w:t xml:space="preserve"> Wie dynamische </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Platform</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> und Spike</w:t></w:r></w:p><w:p $wns><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">2.2.4.8 </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Teleport</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> -</w:t></w:r><w:r><w:t xml:space="preserve"> Diesen dienen zum Teleportieren des Spielers, es gilt folgendes:</w:t></w:r></w:p><w:p $wns><w:r><w:tab/><w:t xml:space="preserve">1. Nach dem Tod wird der Spieler an einen </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Teleporter</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> am Spawn teleportiert</w:t></w:r></w:p><w:p $wns><w:r><w:tab/><w:t xml:space="preserve">2. Der Spieler kann durch </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Teleporter</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> zu einer anderen </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Map</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> transferiert werden mit allen </w:t></w:r><w:r><w:tab/></w:r><w:r><w:tab/><w:t>aktuellen Werten</w:t></w:r></w:p><w:p $wns/>
"@

$xmlFrag4 = $xmlFrag4.Trim()
$rng4.InsertXML($xmlFrag4)

# ---------------------------------------------------------------------------
# 5. After "... dokumentiert werden." insert a new blank, formatted paragraph
#    (bold/bCs/sz24/szCs24/underline) right before "2.3.3 Leistung".
# ---------------------------------------------------------------------------
$idx5 = Find-ParagraphIndexByText $d "dokumentiert werden"
if ($idx5 -lt 0) { throw "Could not find the 'dokumentiert werden' paragraph" }
$p5 = $d.Paragraphs.Item($idx5 + 1)
$start5 = $p5.Range.Start
$rng5 = $d.Range($start5, $start5)
$xmlFrag5 = "<w:p $wns><w:pPr><w:rPr><w:b/><w:bCs/><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/><w:u w:val=`"single`"/></w:rPr></w:pPr></w:p>"
$rng5.InsertXML($xmlFrag5)

Write-Host "All edits applied"
